# Requirements Trace Matrix - progress update on short params
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rows 7-20: unhide (they had been filtered/hidden). Rows 7-12 are now
#     marked "Done" in the Status column (D). Rows 13-20 just become visible. ---
$ws.Range("A7:A20").EntireRow.Hidden = $false

$ws.Range("D7").Value = "D"
$ws.Range("D8").Value = "D"
$ws.Range("D9").Value = "D"
$ws.Range("D10").Value = "D"
$ws.Range("D11").Value = "D"
$ws.Range("D12").Value = "D"

# --- Rows 86-92: mark "Done" while still visible, then hide the rows. ---
$ws.Range("D86").Value = "D"
$ws.Range("D87").Value = "D"
$ws.Range("D88").Value = "D"
$ws.Range("D89").Value = "D"
$ws.Range("D90").Value = "D"
$ws.Range("D91").Value = "D"
$ws.Range("D92").Value = "D"

$ws.Range("A86:A92").EntireRow.Hidden = $true

# --- Rows 146-159 (ShortParameters section): update Status values while
#     temporarily visible, then restore hidden state. ---
$ws.Range("A146:A160").EntireRow.Hidden = $false

$ws.Range("D146").Value = "D"
$ws.Range("D147").Value = "P"
$ws.Range("D149").Value = "D"
$ws.Range("D152").Value = "P"
$ws.Range("D153").Value = "D"
$ws.Range("D154").Value = "P"
$ws.Range("D156").Value = "P"
$ws.Range("D157").Value = "P"
$ws.Range("D158").Value = "P"
$ws.Range("D159").Value = "P"

$ws.Range("A146:A160").EntireRow.Hidden = $true

# --- AutoFilter: switch from PrimeMeridianGeoKey to the Ellipsoid* classes ---
$lo = $ws.ListObjects.Item(1)
$lo.Range.AutoFilter(1, @("requirements_class_EllipsoidGeoKey", "requirements_class_EllipsoidInvFlatteningGeoKey", "requirements_class_EllipsoidSemiMajorAxisGeoKey", "requirements_class_EllipsoidSemiMinorAxisGeoKey"), 7)

# --- Update the active selection ---
$ws.Range("E9").Select()
